$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 308 (pushes the existing rows 308:368 down to 311:371,
# matching the dimension growing from A1:T368 to A1:T371).
$ws.Rows("308:310").Insert()

# Fill the 3 newly inserted rows with the new weekly price-report entries
# (same market/product/category columns as the rest of the sheet; a new date
# group for "Provincia de Limarí" sold in 10-kilo trays).
$commonA = 9
$commonB = "Vega Central Mapocho de Santiago"
$commonC = "Metropolitana"
$commonE = 13
$commonF = "Fruta"
$commonG = 100107
$commonH = "Otros"
$commonI = 100107002
$commonJ = "Chirimoya"
$commonK = "Cultivar IV Región"
$commonD = 45211
$commonQ = "`$/bandeja 10 kilos"
$commonR = "Provincia de Limarí"
$commonT = 10

$newRows = @(
    @{ Row = 308; L = "Especial"; M = 330; N = 27000; O = 27000; P = 27000; S = 2700 },
    @{ Row = 309; L = "Primera";  M = 280; N = 25000; O = 25000; P = 25000; S = 2500 },
    @{ Row = 310; L = "Segunda";  M = 260; N = 20000; O = 20000; P = 20000; S = 2000 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $commonA
    $ws.Cells.Item($row, 2).Value = $commonB
    $ws.Cells.Item($row, 3).Value = $commonC
    $ws.Cells.Item($row, 4).Value = $commonD
    $ws.Cells.Item($row, 5).Value = $commonE
    $ws.Cells.Item($row, 6).Value = $commonF
    $ws.Cells.Item($row, 7).Value = $commonG
    $ws.Cells.Item($row, 8).Value = $commonH
    $ws.Cells.Item($row, 9).Value = $commonI
    $ws.Cells.Item($row, 10).Value = $commonJ
    $ws.Cells.Item($row, 11).Value = $commonK
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $commonQ
    $ws.Cells.Item($row, 18).Value = $commonR
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $commonT
}
